$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update quantity (F) for items with reduced stock; recompute value (G) = Rate(D) * Qty(F) ---
$ws.Range("F18").Value = 4
$ws.Range("G18").Formula = "=D18*F18"
$ws.Range("F23").Value = 1
$ws.Range("G23").Formula = "=D23*F23"
$ws.Range("F42").Value = 161
$ws.Range("G42").Formula = "=D42*F42"
$ws.Range("F67").Value = 34
$ws.Range("G67").Formula = "=D67*F67"
$ws.Range("F108").Value = 11
$ws.Range("G108").Formula = "=D108*F108"
$ws.Range("F195").Value = 48
$ws.Range("G195").Formula = "=D195*F195"
$ws.Range("F197").Value = 31
$ws.Range("G197").Formula = "=D197*F197"
$ws.Range("F198").Value = 57
$ws.Range("G198").Formula = "=D198*F198"
$ws.Range("F202").Value = 52
$ws.Range("G202").Formula = "=D202*F202"
$ws.Range("F203").Value = 207
$ws.Range("G203").Formula = "=D203*F203"
$ws.Range("F204").Value = 44
$ws.Range("G204").Formula = "=D204*F204"
$ws.Range("F205").Value = 102
$ws.Range("G205").Formula = "=D205*F205"
$ws.Range("F206").Value = 5
$ws.Range("G206").Formula = "=D206*F206"
$ws.Range("F208").Value = 16
$ws.Range("G208").Formula = "=D208*F208"
$ws.Range("F211").Value = 237
$ws.Range("G211").Formula = "=D211*F211"
$ws.Range("F212").Value = 2527
$ws.Range("G212").Formula = "=D212*F212"
$ws.Range("F257").Value = 63
$ws.Range("G257").Formula = "=D257*F257"
$ws.Range("F264").Value = 16
$ws.Range("G264").Formula = "=D264*F264"
$ws.Range("F270").Value = 26
$ws.Range("G270").Formula = "=D270*F270"
$ws.Range("F274").Value = 19
$ws.Range("G274").Formula = "=D274*F274"
$ws.Range("F330").Value = 209
$ws.Range("G330").Formula = "=D330*F330"
$ws.Range("F348").Value = 83
$ws.Range("G348").Formula = "=D348*F348"
$ws.Range("F350").Value = 55
$ws.Range("G350").Formula = "=D350*F350"
$ws.Range("F351").Value = 465
$ws.Range("G351").Formula = "=D351*F351"
$ws.Range("F390").Value = 172
$ws.Range("G390").Formula = "=D390*F390"
$ws.Range("F422").Value = 263
$ws.Range("G422").Formula = "=D422*F422"
$ws.Range("F423").Value = 180
$ws.Range("G423").Formula = "=D423*F423"
$ws.Range("F428").Value = 767
$ws.Range("G428").Formula = "=D428*F428"
$ws.Range("F476").Value = 68
$ws.Range("G476").Formula = "=D476*F476"
$ws.Range("F477").Value = 65
$ws.Range("G477").Formula = "=D477*F477"
$ws.Range("F497").Value = 51
$ws.Range("G497").Formula = "=D497*F497"
$ws.Range("F515").Value = 20
$ws.Range("G515").Formula = "=D515*F515"
$ws.Range("F525").Value = 24
$ws.Range("G525").Formula = "=D525*F525"
$ws.Range("F535").Value = 12
$ws.Range("G535").Formula = "=D535*F535"
$ws.Range("F558").Value = 20
$ws.Range("G558").Formula = "=D558*F558"
$ws.Range("F563").Value = 209
$ws.Range("G563").Formula = "=D563*F563"
$ws.Range("F570").Value = 152
$ws.Range("G570").Formula = "=D570*F570"
$ws.Range("F572").Value = 16
$ws.Range("G572").Formula = "=D572*F572"
$ws.Range("F575").Value = 6
$ws.Range("G575").Formula = "=D575*F575"
$ws.Range("F595").Value = 47
$ws.Range("G595").Formula = "=D595*F595"
$ws.Range("F596").Value = 116
$ws.Range("G596").Formula = "=D596*F596"
$ws.Range("F598").Value = 121
$ws.Range("G598").Formula = "=D598*F598"
$ws.Range("F646").Value = 58
$ws.Range("G646").Formula = "=D646*F646"
$ws.Range("F649").Value = 69
$ws.Range("G649").Formula = "=D649*F649"
$ws.Range("F652").Value = 312
$ws.Range("G652").Formula = "=D652*F652"
$ws.Range("F653").Value = 34
$ws.Range("G653").Formula = "=D653*F653"
$ws.Range("F657").Value = 124
$ws.Range("G657").Formula = "=D657*F657"
$ws.Range("F658").Value = 226
$ws.Range("G658").Formula = "=D658*F658"
$ws.Range("F665").Value = 65
$ws.Range("G665").Formula = "=D665*F665"
$ws.Range("F670").Value = 23
$ws.Range("G670").Formula = "=D670*F670"
$ws.Range("F673").Value = 3
$ws.Range("G673").Formula = "=D673*F673"
$ws.Range("F674").Value = 38
$ws.Range("G674").Formula = "=D674*F674"
$ws.Range("F681").Value = 79
$ws.Range("G681").Formula = "=D681*F681"
$ws.Range("F684").Value = 75
$ws.Range("G684").Formula = "=D684*F684"
$ws.Range("F686").Value = 77
$ws.Range("G686").Formula = "=D686*F686"
$ws.Range("F712").Value = 1
$ws.Range("G712").Formula = "=D712*F712"
$ws.Range("F727").Value = 22
$ws.Range("G727").Formula = "=D727*F727"
$ws.Range("F728").Value = 2320
$ws.Range("G728").Formula = "=D728*F728"
$ws.Range("F729").Value = 270
$ws.Range("G729").Formula = "=D729*F729"
$ws.Range("F730").Value = 367
$ws.Range("G730").Formula = "=D730*F730"
$ws.Range("F731").Value = 62
$ws.Range("G731").Formula = "=D731*F731"
$ws.Range("F732").Value = 138
$ws.Range("G732").Formula = "=D732*F732"
$ws.Range("F734").Value = 125
$ws.Range("G734").Formula = "=D734*F734"

# --- Row 152/153: swap Code (B) values between two identical product rows ---
$tmpB152 = $ws.Range("B152").Value2
$tmpB153 = $ws.Range("B153").Value2
$ws.Range("B152").Value = $tmpB153
$ws.Range("B153").Value = $tmpB152

# --- Row 258/259: swap entire item data (Code, Rate, MRP, Qty, Value) between two rows ---
$tmpB258 = $ws.Range("B258").Value2
$tmpB259 = $ws.Range("B259").Value2
$tmpD258 = $ws.Range("D258").Value2
$tmpD259 = $ws.Range("D259").Value2
$tmpE258 = $ws.Range("E258").Value2
$tmpE259 = $ws.Range("E259").Value2
$tmpF258 = $ws.Range("F258").Value2
$tmpF259 = $ws.Range("F259").Value2
$tmpG258 = $ws.Range("G258").Value2
$tmpG259 = $ws.Range("G259").Value2
$ws.Range("B258").Value = $tmpB259
$ws.Range("B259").Value = $tmpB258
$ws.Range("D258").Value = $tmpD259
$ws.Range("D259").Value = $tmpD258
$ws.Range("E258").Value = $tmpE259
$ws.Range("E259").Value = $tmpE258
$ws.Range("F258").Value = $tmpF259
$ws.Range("F259").Value = $tmpF258
$ws.Range("G258").Value = $tmpG259
$ws.Range("G259").Value = $tmpG258

# --- Recompute all "Sub Total:" rows as SUM of their item Value (G) range ---
$ws.Range("B7").Formula = "=SUM(G5:G6)"
$ws.Range("B12").Formula = "=SUM(G9:G11)"
$ws.Range("B25").Formula = "=SUM(G14:G24)"
$ws.Range("B54").Formula = "=SUM(G27:G53)"
$ws.Range("B57").Formula = "=SUM(G56:G56)"
$ws.Range("B69").Formula = "=SUM(G59:G68)"
$ws.Range("B72").Formula = "=SUM(G71:G71)"
$ws.Range("B81").Formula = "=SUM(G74:G80)"
$ws.Range("B116").Formula = "=SUM(G83:G115)"
$ws.Range("B119").Formula = "=SUM(G118:G118)"
$ws.Range("B123").Formula = "=SUM(G121:G122)"
$ws.Range("B131").Formula = "=SUM(G125:G130)"
$ws.Range("B136").Formula = "=SUM(G133:G135)"
$ws.Range("B147").Formula = "=SUM(G138:G146)"
$ws.Range("B154").Formula = "=SUM(G149:G153)"
$ws.Range("B161").Formula = "=SUM(G156:G160)"
$ws.Range("B176").Formula = "=SUM(G163:G175)"
$ws.Range("B187").Formula = "=SUM(G178:G186)"
$ws.Range("B193").Formula = "=SUM(G189:G192)"
$ws.Range("B209").Formula = "=SUM(G195:G208)"
$ws.Range("B217").Formula = "=SUM(G211:G216)"
$ws.Range("B226").Formula = "=SUM(G219:G225)"
$ws.Range("B231").Formula = "=SUM(G228:G230)"
$ws.Range("B254").Formula = "=SUM(G233:G253)"
$ws.Range("B313").Formula = "=SUM(G256:G312)"
$ws.Range("B321").Formula = "=SUM(G315:G320)"
$ws.Range("B327").Formula = "=SUM(G323:G326)"
$ws.Range("B336").Formula = "=SUM(G329:G335)"
$ws.Range("B346").Formula = "=SUM(G338:G345)"
$ws.Range("B353").Formula = "=SUM(G348:G352)"
$ws.Range("B360").Formula = "=SUM(G355:G359)"
$ws.Range("B369").Formula = "=SUM(G362:G368)"
$ws.Range("B372").Formula = "=SUM(G371:G371)"
$ws.Range("B385").Formula = "=SUM(G374:G384)"
$ws.Range("B409").Formula = "=SUM(G387:G408)"
$ws.Range("B412").Formula = "=SUM(G411:G411)"
$ws.Range("B415").Formula = "=SUM(G414:G414)"
$ws.Range("B435").Formula = "=SUM(G417:G434)"
$ws.Range("B457").Formula = "=SUM(G437:G456)"
$ws.Range("B463").Formula = "=SUM(G459:G462)"
$ws.Range("B467").Formula = "=SUM(G465:G466)"
$ws.Range("B473").Formula = "=SUM(G469:G472)"
$ws.Range("B479").Formula = "=SUM(G475:G478)"
$ws.Range("B489").Formula = "=SUM(G481:G488)"
$ws.Range("B493").Formula = "=SUM(G491:G492)"
$ws.Range("B507").Formula = "=SUM(G495:G506)"
$ws.Range("B526").Formula = "=SUM(G509:G525)"
$ws.Range("B540").Formula = "=SUM(G528:G539)"
$ws.Range("B552").Formula = "=SUM(G542:G551)"
$ws.Range("B576").Formula = "=SUM(G554:G575)"
$ws.Range("B584").Formula = "=SUM(G578:G583)"
$ws.Range("B592").Formula = "=SUM(G586:G591)"
$ws.Range("B599").Formula = "=SUM(G594:G598)"
$ws.Range("B609").Formula = "=SUM(G601:G608)"
$ws.Range("B612").Formula = "=SUM(G611:G611)"
$ws.Range("B618").Formula = "=SUM(G614:G617)"
$ws.Range("B622").Formula = "=SUM(G620:G621)"
$ws.Range("B626").Formula = "=SUM(G624:G625)"
$ws.Range("B635").Formula = "=SUM(G628:G634)"
$ws.Range("B638").Formula = "=SUM(G637:G637)"
$ws.Range("B644").Formula = "=SUM(G640:G643)"
$ws.Range("B660").Formula = "=SUM(G646:G659)"
$ws.Range("B690").Formula = "=SUM(G662:G689)"
$ws.Range("B696").Formula = "=SUM(G692:G695)"
$ws.Range("B705").Formula = "=SUM(G698:G704)"
$ws.Range("B721").Formula = "=SUM(G707:G720)"
$ws.Range("B725").Formula = "=SUM(G723:G724)"
$ws.Range("B736").Formula = "=SUM(G727:G735)"
$ws.Range("B740").Formula = "=SUM(G738:G739)"

# --- Recompute grand Sub Total (741) as SUM of all company subtotal rows ---
$ws.Range("B741").Formula = "=SUM(B7,B12,B25,B54,B57,B69,B72,B81,B116,B119,B123,B131,B136,B147,B154,B161,B176,B187,B193,B209,B217,B226,B231,B254,B313,B321,B327,B336,B346,B353,B360,B369,B372,B385,B409,B412,B415,B435,B457,B463,B467,B473,B479,B489,B493,B507,B526,B540,B552,B576,B584,B592,B599,B609,B612,B618,B622,B626,B635,B638,B644,B660,B690,B696,B705,B721,B725,B736,B740)"

# --- Grand Total (742) mirrors Sub Total (741) ---
$ws.Range("B742").Formula = "=B741"
